$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2189.2
$ws.Range("J17").Value = 2189.2
$ws.Range("L17").Value = 6567.599999999999
$ws.Range("N17").Value = -6903.599999999999
$ws.Range("H44").Value = 12500
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10924
$ws.Range("H52").Value = 9750
$ws.Range("I52").Value = 9750
$ws.Range("K52").Value = 29250
$ws.Range("M52").Value = -29090
$ws.Range("H76").Value = 5380.154
$ws.Range("I76").Value = 6242.5713
$ws.Range("K76").Value = 6242.5713
$ws.Range("M76").Value = -5927.5713
$ws.Range("H79").Value = 5380.154
$ws.Range("I79").Value = 6242.5713
$ws.Range("K79").Value = 6242.5713
$ws.Range("M79").Value = -5150.5713
$ws.Range("H86").Value = 4713.75
$ws.Range("I86").Value = 4670.2856
$ws.Range("K86").Value = 4670.2856
$ws.Range("M86").Value = -3547.2856
$ws.Range("H89").Value = 4713.75
$ws.Range("I89").Value = 4670.2856
$ws.Range("K89").Value = 23351.428
$ws.Range("M89").Value = -17735.428
$ws.Range("H100").Value = 4646.1304
$ws.Range("I100").Value = 2268.0908
$ws.Range("J100").Value = 6826
$ws.Range("K100").Value = 2268.0908
$ws.Range("L100").Value = 6826
$ws.Range("M100").Value = -1727.0908
$ws.Range("N100").Value = -7908
$ws.Range("H107").Value = 594.8
$ws.Range("I107").Value = 659.6667
$ws.Range("K107").Value = 659.6667
$ws.Range("M107").Value = 1260.3333
$ws.Range("H138").Value = 2831.8357
$ws.Range("I138").Value = 1344.4
$ws.Range("J138").Value = 3606.5417
$ws.Range("K138").Value = 4033.2
$ws.Range("L138").Value = 10819.6251
$ws.Range("M138").Value = 1106.8
$ws.Range("N138").Value = -21099.6251

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9240.569
$ws.Range("I32").Value = 8439.119000000001
$ws.Range("J32").Value = 19980
$ws.Range("K32").Value = 8439.119000000001
$ws.Range("L32").Value = 19980
$ws.Range("M32").Value = -8152.119000000001
$ws.Range("N32").Value = -20554
$ws.Range("H61").Value = 3377.1177
$ws.Range("I61").Value = 1142.2
$ws.Range("J61").Value = 4308.3335
$ws.Range("K61").Value = 1142.2
$ws.Range("L61").Value = 4308.3335
$ws.Range("M61").Value = -930.2
$ws.Range("N61").Value = -4732.3335
$ws.Range("H74").Value = 1873.9546
$ws.Range("I74").Value = 1930
$ws.Range("K74").Value = 1930
$ws.Range("M74").Value = -1056
$ws.Range("H77").Value = 1873.9546
$ws.Range("I77").Value = 1930
$ws.Range("K77").Value = 9650
$ws.Range("M77").Value = -5282
$ws.Range("H110").Value = 3818.1333
$ws.Range("I110").Value = 3733.7144
$ws.Range("K110").Value = 3733.7144
$ws.Range("M110").Value = -1688.7144
$ws.Range("H132").Value = 3615.575
$ws.Range("I132").Value = 3394.6333
$ws.Range("K132").Value = 10183.8999
$ws.Range("M132").Value = -7653.8999
$ws.Range("H133").Value = 68490.2
$ws.Range("I133").Value = 63206.168
$ws.Range("J133").Value = 69868.64999999999
$ws.Range("K133").Value = 63206.168
$ws.Range("L133").Value = 69868.64999999999
$ws.Range("M133").Value = -60676.168
$ws.Range("N133").Value = -74928.64999999999
$ws.Range("H136").Value = 3377.1177
$ws.Range("I136").Value = 1142.2
$ws.Range("J136").Value = 4308.3335
$ws.Range("K136").Value = 3426.6
$ws.Range("L136").Value = 12925.0005
$ws.Range("M136").Value = -876.6000000000004
$ws.Range("N136").Value = -18025.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3847.111
$ws.Range("I86").Value = 1565.6666
$ws.Range("K86").Value = 1565.6666
$ws.Range("M86").Value = -442.6666
$ws.Range("H89").Value = 3847.111
$ws.Range("I89").Value = 1565.6666
$ws.Range("K89").Value = 7828.333000000001
$ws.Range("M89").Value = -2212.333000000001
$ws.Range("H105").Value = 2403.2104
$ws.Range("I105").Value = 2550.2
$ws.Range("J105").Value = 1852
$ws.Range("K105").Value = 2550.2
$ws.Range("L105").Value = 1852
$ws.Range("M105").Value = -803.1999999999998
$ws.Range("N105").Value = -5346
$ws.Range("H134").Value = 2519.6445
$ws.Range("I134").Value = 1264.3823
$ws.Range("K134").Value = 3793.1469
$ws.Range("M134").Value = -1258.1469

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3909.2693
$ws.Range("I132").Value = 2936.9546
$ws.Range("K132").Value = 8810.863799999999
$ws.Range("M132").Value = -6280.863799999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 4000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 4000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 12000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -12962
$ws.Range("H53").Value = 4000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 4000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 12000
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -12962

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2083.0264
$ws.Range("I132").Value = 2202.3
$ws.Range("J132").Value = 1635.75
$ws.Range("K132").Value = 6606.900000000001
$ws.Range("L132").Value = 4907.25
$ws.Range("M132").Value = -4076.900000000001
$ws.Range("N132").Value = -9967.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9026.727999999999
$ws.Range("I46").Value = 2250
$ws.Range("J46").Value = 9961.448
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 9961.448
$ws.Range("M46").Value = -2062
$ws.Range("N46").Value = -10337.448
$ws.Range("H55").Value = 5090.8
$ws.Range("J55").Value = 1734.6666
$ws.Range("L55").Value = 1734.6666
$ws.Range("N55").Value = -2080.6666
$ws.Range("H82").Value = 3366
$ws.Range("I82").Value = 2034.2778
$ws.Range("J82").Value = 6790.4287
$ws.Range("K82").Value = 2034.2778
$ws.Range("L82").Value = 6790.4287
$ws.Range("M82").Value = -1673.2778
$ws.Range("N82").Value = -7512.4287
$ws.Range("H85").Value = 3366
$ws.Range("I85").Value = 2034.2778
$ws.Range("J85").Value = 6790.4287
$ws.Range("K85").Value = 2034.2778
$ws.Range("L85").Value = 6790.4287
$ws.Range("M85").Value = -786.2778000000001
$ws.Range("N85").Value = -9286.4287
$ws.Range("H132").Value = 4005.4211
$ws.Range("I132").Value = 3422.8147
$ws.Range("K132").Value = 10268.4441
$ws.Range("M132").Value = -7738.444100000001
$ws.Range("H136").Value = 4423.625
$ws.Range("I136").Value = 2283
$ws.Range("J136").Value = 6360.381
$ws.Range("K136").Value = 6849
$ws.Range("L136").Value = 19081.143
$ws.Range("M136").Value = -4299
$ws.Range("N136").Value = -24181.143

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17995
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H113").Value = 429.57693
$ws.Range("I113").Value = 400.17648
$ws.Range("J113").Value = 485.1111
$ws.Range("K113").Value = 1200.52944
$ws.Range("L113").Value = 1455.3333
$ws.Range("M113").Value = 969.47056
$ws.Range("N113").Value = -5795.3333
$ws.Range("H136").Value = 2176.2195
$ws.Range("I136").Value = 1074.2258
$ws.Range("K136").Value = 3222.6774
$ws.Range("M136").Value = -672.6773999999996
